$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7.76996414094062
$ws.Cells.Item(2, 3).Value = 5.950183189536349
$ws.Cells.Item(2, 5).Value = 13.15868052146605
$ws.Cells.Item(2, 6).Value = 16.86991607391245
$ws.Cells.Item(2, 7).Value = 22.48454416139465
$ws.Cells.Item(2, 8).Value = 12.89568284388557
$ws.Cells.Item(2, 9).Value = 18.51254955914431
$ws.Cells.Item(2, 11).Value = 8.756972058867465
$ws.Cells.Item(2, 13).Value = 12.92701301089317
$ws.Cells.Item(2, 15).Value = 18.7619200971055
$ws.Cells.Item(3, 2).Value = 7.373936293603952
$ws.Cells.Item(3, 3).Value = 5.840722907915868
$ws.Cells.Item(3, 5).Value = 12.96737368147768
$ws.Cells.Item(3, 6).Value = 15.89584955866815
$ws.Cells.Item(3, 7).Value = 22.65360064131114
$ws.Cells.Item(3, 8).Value = 12.95186815686374
$ws.Cells.Item(3, 9).Value = 18.62343492872452
$ws.Cells.Item(3, 11).Value = 8.514745323644695
$ws.Cells.Item(3, 13).Value = 12.72954518438232
$ws.Cells.Item(3, 15).Value = 18.86783593566925
$ws.Cells.Item(4, 2).Value = 7.11819658360023
$ws.Cells.Item(4, 3).Value = 5.772084483911598
$ws.Cells.Item(4, 5).Value = 12.85344494733653
$ws.Cells.Item(4, 6).Value = 15.26997757108489
$ws.Cells.Item(4, 7).Value = 22.76579128763034
$ws.Cells.Item(4, 8).Value = 12.98838321378037
$ws.Cells.Item(4, 9).Value = 18.69521970370205
$ws.Cells.Item(4, 11).Value = 8.361243886109346
$ws.Cells.Item(4, 13).Value = 12.60941274055197
$ws.Cells.Item(4, 15).Value = 18.9370219302552
$ws.Cells.Item(5, 2).Value = 7.010886177909661
$ws.Cells.Item(5, 3).Value = 5.743779119643774
$ws.Cells.Item(5, 5).Value = 12.80797213612557
$ws.Cells.Item(5, 6).Value = 15.008197319934
$ws.Cells.Item(5, 7).Value = 22.81360913049919
$ws.Cells.Item(5, 8).Value = 13.00377128297473
$ws.Cells.Item(5, 9).Value = 18.72540468092885
$ws.Cells.Item(5, 11).Value = 8.29755130468175
$ws.Cells.Item(5, 13).Value = 12.56079910879255
$ws.Cells.Item(5, 15).Value = 18.96625986032389
$ws.Cells.Item(6, 2).Value = 6.992882730296933
$ws.Cells.Item(6, 3).Value = 5.73905956641196
$ws.Cells.Item(6, 5).Value = 12.80048089857949
$ws.Cells.Item(6, 6).Value = 14.96433081551589
$ws.Cells.Item(6, 7).Value = 22.82167573743813
$ws.Cells.Item(6, 8).Value = 13.00635715566068
$ws.Cells.Item(6, 9).Value = 18.73047321527299
$ws.Cells.Item(6, 11).Value = 8.286908247934145
$ws.Cells.Item(6, 13).Value = 12.55274925363583
$ws.Cells.Item(6, 15).Value = 18.97117784855183
$ws.Cells.Item(7, 2).Value = 7.11676178717189
$ws.Cells.Item(7, 3).Value = 5.771704069406613
$ws.Cells.Item(7, 5).Value = 12.85282773745683
$ws.Cells.Item(7, 6).Value = 15.26647399323133
$ws.Cells.Item(7, 7).Value = 22.76642768955178
$ws.Cells.Item(7, 8).Value = 12.98858868533195
$ws.Cells.Item(7, 9).Value = 18.69562301278586
$ws.Cells.Item(7, 11).Value = 8.360389435631186
$ws.Cells.Item(7, 13).Value = 12.60875565732701
$ws.Cells.Item(7, 15).Value = 18.93741201613033
$ws.Cells.Item(8, 2).Value = 7.636061034266178
$ws.Cells.Item(8, 3).Value = 5.912752506539266
$ws.Cells.Item(8, 5).Value = 13.09202417443119
$ws.Cells.Item(8, 6).Value = 16.5399640634477
$ws.Cells.Item(8, 7).Value = 22.5410874039119
$ws.Cells.Item(8, 8).Value = 12.91463731246003
$ws.Cells.Item(8, 9).Value = 18.55001566895622
$ws.Cells.Item(8, 11).Value = 8.674477926773797
$ws.Cells.Item(8, 13).Value = 12.85873177604054
$ws.Cells.Item(8, 15).Value = 18.79757768966983
$ws.Cells.Item(9, 2).Value = 8.552233314612792
$ws.Cells.Item(9, 3).Value = 6.176971004598331
$ws.Cells.Item(9, 5).Value = 13.58616981982573
$ws.Cells.Item(9, 6).Value = 19.00274580682531
$ws.Cells.Item(9, 7).Value = 22.16622717623811
$ws.Cells.Item(9, 8).Value = 12.78559196656377
$ws.Cells.Item(9, 9).Value = 18.29377350764356
$ws.Cells.Item(9, 11).Value = 9.25010760625025
$ws.Cells.Item(9, 13).Value = 13.3550817739148
$ws.Cells.Item(9, 15).Value = 18.55633039776461
$ws.Cells.Item(10, 2).Value = 9.160588682987481
$ws.Cells.Item(10, 3).Value = 6.362239398663512
$ws.Cells.Item(10, 5).Value = 13.96046702321408
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 21.9323273320692
$ws.Cells.Item(10, 8).Value = 12.70047608958126
$ws.Cells.Item(10, 9).Value = 18.12327180222737
$ws.Cells.Item(10, 11).Value = 9.645548618939513
$ws.Cells.Item(10, 13).Value = 13.71993905098042
$ws.Cells.Item(10, 15).Value = 18.39920265070506
$ws.Cells.Item(11, 2).Value = 9.422927891957205
$ws.Cells.Item(11, 3).Value = 6.444347140043043
$ws.Cells.Item(11, 5).Value = 14.13232620168121
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 21.83507543489763
$ws.Cells.Item(11, 8).Value = 12.66385056525596
$ws.Cells.Item(11, 9).Value = 18.04954211043944
$ws.Cells.Item(11, 11).Value = 9.818955828545064
$ws.Cells.Item(11, 13).Value = 13.8852315728951
$ws.Cells.Item(11, 15).Value = 18.33209470947934
$ws.Cells.Item(12, 2).Value = 9.520175910087124
$ws.Cells.Item(12, 3).Value = 6.475107350644313
$ws.Cells.Item(12, 5).Value = 14.19756267831148
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 21.79957652146926
$ws.Cells.Item(12, 8).Value = 12.65028195692421
$ws.Cells.Item(12, 9).Value = 18.02217230596693
$ws.Cells.Item(12, 11).Value = 9.883650655125193
$ws.Cells.Item(12, 13).Value = 13.94766661398791
$ws.Cells.Item(12, 15).Value = 18.30731183689906
$ws.Cells.Item(13, 2).Value = 9.499325254793375
$ws.Cells.Item(13, 3).Value = 6.468497667590858
$ws.Cells.Item(13, 5).Value = 14.18350696708443
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 21.80716257265713
$ws.Cells.Item(13, 8).Value = 12.65319083344712
$ws.Cells.Item(13, 9).Value = 18.02804243673809
$ws.Cells.Item(13, 11).Value = 9.869761214504555
$ws.Cells.Item(13, 13).Value = 13.93422804559313
$ws.Cells.Item(13, 15).Value = 18.31262126017032
$ws.Cells.Item(14, 2).Value = 9.430970626735686
$ws.Cells.Item(14, 3).Value = 6.446884568597581
$ws.Cells.Item(14, 5).Value = 14.13769050001474
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 21.8321282255748
$ws.Cells.Item(14, 8).Value = 12.66272824403717
$ws.Cells.Item(14, 9).Value = 18.04727936264812
$ws.Cells.Item(14, 11).Value = 9.824297973368141
$ws.Cells.Item(14, 13).Value = 13.8903716336046
$ws.Cells.Item(14, 15).Value = 18.33004318650421
$ws.Cells.Item(15, 2).Value = 9.388828137007597
$ws.Cells.Item(15, 3).Value = 6.433602088633263
$ws.Cells.Item(15, 5).Value = 14.10964496210848
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 21.84759376088162
$ws.Cells.Item(15, 8).Value = 12.66860932827717
$ws.Cells.Item(15, 9).Value = 18.05913413463734
$ws.Cells.Item(15, 11).Value = 9.796322926936845
$ws.Cells.Item(15, 13).Value = 13.86348605694195
$ws.Cells.Item(15, 15).Value = 18.34079662343525
$ws.Cells.Item(16, 2).Value = 9.14315227126394
$ws.Cells.Item(16, 3).Value = 6.356828032138553
$ws.Cells.Item(16, 5).Value = 13.94926146516917
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 21.9388681625812
$ws.Cells.Item(16, 8).Value = 12.7029117538334
$ws.Cells.Item(16, 9).Value = 18.12816723639999
$ws.Cells.Item(16, 11).Value = 9.634082230548318
$ws.Cells.Item(16, 13).Value = 13.70911792369346
$ws.Cells.Item(16, 15).Value = 18.40367633795826
$ws.Cells.Item(17, 2).Value = 8.988730772852616
$ws.Cells.Item(17, 3).Value = 6.309158823957786
$ws.Cells.Item(17, 5).Value = 13.8512272001192
$ws.Cells.Item(17, 6).Value = 20.20408069597325
$ws.Cells.Item(17, 7).Value = 21.99721449406832
$ws.Cells.Item(17, 8).Value = 12.72449120773818
$ws.Cells.Item(17, 9).Value = 18.17149760256325
$ws.Cells.Item(17, 11).Value = 9.532864021118908
$ws.Cells.Item(17, 13).Value = 13.61420020285045
$ws.Cells.Item(17, 15).Value = 18.4433710298763
$ws.Cells.Item(18, 2).Value = 8.898556760603169
$ws.Cells.Item(18, 3).Value = 6.281537448566782
$ws.Cells.Item(18, 5).Value = 13.79499428987915
$ws.Cells.Item(18, 6).Value = 19.95656407809801
$ws.Cells.Item(18, 7).Value = 22.03163427753023
$ws.Cells.Item(18, 8).Value = 12.73710027224558
$ws.Cells.Item(18, 9).Value = 18.19678086611933
$ws.Cells.Item(18, 11).Value = 9.47403802365003
$ws.Cells.Item(18, 13).Value = 13.55954519444369
$ws.Cells.Item(18, 15).Value = 18.46661365118889
$ws.Cells.Item(19, 2).Value = 8.867793474428266
$ws.Cells.Item(19, 3).Value = 6.272151034490753
$ws.Cells.Item(19, 5).Value = 13.77598341030212
$ws.Cells.Item(19, 6).Value = 19.87204792380568
$ws.Cells.Item(19, 7).Value = 22.04343564175089
$ws.Cells.Item(19, 8).Value = 12.74140335672402
$ws.Cells.Item(19, 9).Value = 18.20540334280021
$ws.Cells.Item(19, 11).Value = 9.454017382078861
$ws.Cells.Item(19, 13).Value = 13.54103134360045
$ws.Cells.Item(19, 15).Value = 18.47455381233137
$ws.Cells.Item(20, 2).Value = 9.005309608558713
$ws.Cells.Item(20, 3).Value = 6.314254482350028
$ws.Cells.Item(20, 5).Value = 13.86164771904513
$ws.Cells.Item(20, 6).Value = 20.24955283636154
$ws.Cells.Item(20, 7).Value = 21.99091427732561
$ws.Cells.Item(20, 8).Value = 12.72217363878069
$ws.Cells.Item(20, 9).Value = 18.16684767928445
$ws.Cells.Item(20, 11).Value = 9.543702087645739
$ws.Cells.Item(20, 13).Value = 13.62431109105808
$ws.Cells.Item(20, 15).Value = 18.43910288867076
$ws.Cells.Item(21, 2).Value = 9.451105028204559
$ws.Cells.Item(21, 3).Value = 6.453242028148547
$ws.Cells.Item(21, 5).Value = 14.15114420359785
$ws.Cells.Item(21, 6).Value = 21.46857628470577
$ws.Cells.Item(21, 7).Value = 21.82475905931151
$ws.Cells.Item(21, 8).Value = 12.65991871837002
$ws.Cells.Item(21, 9).Value = 18.04161409177131
$ws.Cells.Item(21, 11).Value = 9.837678263508955
$ws.Cells.Item(21, 13).Value = 13.90325806643816
$ws.Cells.Item(21, 15).Value = 18.32490885295937
$ws.Cells.Item(22, 2).Value = 9.730244743418075
$ws.Cells.Item(22, 3).Value = 6.542133710171615
$ws.Cells.Item(22, 5).Value = 14.34123091252304
$ws.Cells.Item(22, 6).Value = 22.22866616901552
$ws.Cells.Item(22, 7).Value = 21.72391506053783
$ws.Cells.Item(22, 8).Value = 12.6209839470609
$ws.Cells.Item(22, 9).Value = 17.96297223980613
$ws.Cells.Item(22, 11).Value = 10.02413538431492
$ws.Cells.Item(22, 13).Value = 14.08462115241313
$ws.Cells.Item(22, 15).Value = 18.25394587749267
$ws.Cells.Item(23, 2).Value = 9.582386365656667
$ws.Cells.Item(23, 3).Value = 6.494874739307726
$ws.Cells.Item(23, 5).Value = 14.23972032690105
$ws.Cells.Item(23, 6).Value = 21.82633154458858
$ws.Cells.Item(23, 7).Value = 21.77702430357407
$ws.Cells.Item(23, 8).Value = 12.64160396309068
$ws.Cells.Item(23, 9).Value = 18.00465190737501
$ws.Cells.Item(23, 11).Value = 9.925150601103178
$ws.Cells.Item(23, 13).Value = 13.98792978286422
$ws.Cells.Item(23, 15).Value = 18.29148404701339
$ws.Cells.Item(24, 2).Value = 8.997818660403999
$ws.Cells.Item(24, 3).Value = 6.311951405580269
$ws.Cells.Item(24, 5).Value = 13.8569361973612
$ws.Cells.Item(24, 6).Value = 20.22900810905287
$ws.Cells.Item(24, 7).Value = 21.99375987850619
$ws.Cells.Item(24, 8).Value = 12.7232207801161
$ws.Cells.Item(24, 9).Value = 18.16894875057023
$ws.Cells.Item(24, 11).Value = 9.538804170511153
$ws.Cells.Item(24, 13).Value = 13.61974022039879
$ws.Cells.Item(24, 15).Value = 18.44103120207527
$ws.Cells.Item(25, 2).Value = 8.315611873147366
$ws.Cells.Item(25, 3).Value = 6.106957464342857
$ws.Cells.Item(25, 5).Value = 13.45022094327532
$ws.Cells.Item(25, 6).Value = 18.34778573295695
$ws.Cells.Item(25, 7).Value = 22.26039454635547
$ws.Cells.Item(25, 8).Value = 12.81879674118956
$ws.Cells.Item(25, 9).Value = 18.35996811206248
$ws.Cells.Item(25, 11).Value = 9.099016868962927
$ws.Cells.Item(25, 13).Value = 13.22051916333757
$ws.Cells.Item(25, 15).Value = 18.31262126017032
